$d = $word.ActiveDocument

# --- "Previous" row: rename to "Previously" and bump the row height 607 -> 624 twips ---
$findRange1 = $d.Content.Duplicate
$findRange1.Find.Execute("Previous", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$row1 = $findRange1.Cells.Item(1).Row
$row1.Height = 31.2            # 624 twips == 31.2 points (Word reports/accepts Height in points)
$findRange1.Find.Execute("Previous", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "Previously", 2)

# --- "Current" row: rename to "Currently" and bump the row height 607 -> 624 twips ---
$findRange2 = $d.Content.Duplicate
$findRange2.Find.Execute("Current", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$row2 = $findRange2.Cells.Item(1).Row
$row2.Height = 31.2             # 624 twips == 31.2 points
$findRange2.Find.Execute("Current", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "Currently", 2)
